$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert "Instagram password needed to be lowercase":
# C3 held the lowercase "peoplespaceoc" Instagram password; restore the
# original capitalised "Peoplespaceoc" used for the other passwords.
$ws.Range("C3").Value = "Peoplespaceoc"
